$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.976.39"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").Value = "2.752.71"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'351.45"
$ws.Range("E5").Value = "  -1.78%  "
$ws.Range("D6").Value = "'107.11"
$ws.Range("E6").Value = "  -1.98%  "
$ws.Range("E7").Value = "  -2.51%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.80%  "
$ws.Range("D10").Value = "'39.01"
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("E11").Value = "  +3.54%  "
$ws.Range("D12").Value = "'0.0831"
$ws.Range("E12").Value = "  -2.31%  "
$ws.Range("D13").Value = "'19.60"
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("E14").Value = "  -2.19%  "
$ws.Range("D15").Value = "3.183.03"
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").Value = "2.740.30"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").Value = "'0.923"
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").Value = "50.895.28"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").Value = "'7.59"
$ws.Range("E19").Value = "  +2.81%  "
$ws.Range("E20").Value = "  -2.67%  "
$ws.Range("E21").Value = "  -1.30%  "
$ws.Range("E22").Value = "  -2.62%  "
$ws.Range("D23").Value = "'69.10"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("D24").Value = "'262.71"
$ws.Range("E24").Value = "  -3.80%  "
$ws.Range("D25").Value = "'2.70"
$ws.Range("E25").Value = "  -1.91%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").Value = "'25.78"
$ws.Range("E27").Value = "  -2.51%  "
$ws.Range("E28").Value = "  +13.11%  "
$ws.Range("E29").Value = "  +2.31%  "
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("D31").Value = "'51.40"
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("D32").Value = "'34.30"
$ws.Range("E32").Value = "  +1.17%  "
$ws.Range("E34").Value = "  -7.45%  "
$ws.Range("D35").Value = "'5.28"
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "'18.20"
$ws.Range("E38").Value = "  +1.13%  "
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("E40").Value = "  -2.84%  "
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("D42").Value = "'2.47"
$ws.Range("E42").Value = "  -5.01%  "
$ws.Range("D43").Value = "'120.63"
$ws.Range("E43").Value = "  -3.54%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'2.19"
$ws.Range("E44").Value = "  -2.67%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'21.79"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("D46").Value = "2.085.01"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("B49").Value = "SEI"
$ws.Range("C49").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D49").Value = "'0.906"
$ws.Range("E49").Value = "  -3.03%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "'5.41"
$ws.Range("E50").Value = "  -5.00%  "
$ws.Range("E51").Value = "  +4.61%  "
